# Updated cryptos list on Tue Aug 20 17:10:39 UTC 2024 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row with
# the latest scraped figures. Values that look like a plain decimal number
# (e.g. "561.66") are written with a leading apostrophe so Excel keeps them
# as literal text -- matching the existing cell contents, which use '.' as
# a thousands separator in some rows (e.g. "58.811.51") and must stay text
# throughout the column instead of being auto-coerced to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new Price, new Volume(1h))
$updates = @(
    @{ Row = 2;  Price = "58.811.51";  Volume = "  +0.40%  " },
    @{ Row = 3;  Price = "2.570.74";   Volume = "  -0.40%  " },
    @{ Row = 4;  Price = $null;        Volume = "  -0.08%  " },
    @{ Row = 5;  Price = "561.66";     Volume = "  +2.04%  " },
    @{ Row = 6;  Price = "142.89";     Volume = "  -1.08%  " },
    @{ Row = 7;  Price = $null;        Volume = "  +0.00%  " },
    @{ Row = 8;  Price = "0.599";      Volume = "  +1.80%  " },
    @{ Row = 9;  Price = "2.572.55";   Volume = "  -0.51%  " },
    @{ Row = 10; Price = "6.67";       Volume = "  -2.57%  " },
    @{ Row = 11; Price = $null;        Volume = "  +2.82%  " },
    @{ Row = 12; Price = $null;        Volume = "  +8.31%  " },
    @{ Row = 13; Price = $null;        Volume = "  +2.35%  " },
    @{ Row = 14; Price = "3.020.88";   Volume = "  -0.45%  " },
    @{ Row = 15; Price = "58.884.45";  Volume = "  +0.63%  " },
    @{ Row = 16; Price = "22.01";      Volume = "  +6.54%  " },
    @{ Row = 17; Price = $null;        Volume = "  +4.39%  " },
    @{ Row = 18; Price = "2.575.85";   Volume = "  -0.50%  " },
    @{ Row = 19; Price = "4.49";       Volume = "  +1.30%  " },
    @{ Row = 20; Price = "335.38";     Volume = "  +0.16%  " },
    @{ Row = 21; Price = "10.16";      Volume = "  +1.24%  " },
    @{ Row = 22; Price = "6.15";       Volume = "  +0.87%  " },
    @{ Row = 23; Price = $null;        Volume = "  -0.16%  " },
    @{ Row = 24; Price = "63.59";      Volume = "  -4.47%  " },
    @{ Row = 25; Price = "0.451";      Volume = "  +6.59%  " },
    @{ Row = 26; Price = "0.999";      Volume = "  +0.15%  " },
    @{ Row = 27; Price = "0.161";      Volume = "  +1.97%  " },
    @{ Row = 28; Price = "7.23";       Volume = "  +2.66%  " },
    @{ Row = 29; Price = "0.0₃0778";   Volume = "  +5.08%  " },
    @{ Row = 30; Price = "0.999";      Volume = "  +0.00%  " },
    @{ Row = 31; Price = $null;        Volume = "  +0.58%  " },
    @{ Row = 32; Price = "158.29";     Volume = "  +1.92%  " },
    @{ Row = 33; Price = $null;        Volume = "  +2.40%  " },
    @{ Row = 34; Price = $null;        Volume = "  +0.63%  " },
    @{ Row = 35; Price = "3.99";       Volume = "  +2.26%  " },
    @{ Row = 36; Price = "0.876";      Volume = "  +2.55%  " },
    @{ Row = 37; Price = "0.879";      Volume = "  +7.44%  " },
    @{ Row = 38; Price = "1.13";       Volume = "  +2.05%  " },
    @{ Row = 39; Price = "36.71";      Volume = "  -1.37%  " },
    @{ Row = 40; Price = $null;        Volume = "  +3.63%  " },
    @{ Row = 41; Price = "290.94";     Volume = "  +4.25%  " },
    @{ Row = 42; Price = "3.64";       Volume = "  +1.49%  " },
    @{ Row = 43; Price = $null;        Volume = "  +0.12%  " },
    @{ Row = 44; Price = "0.0971";     Volume = "  +2.57%  " },
    @{ Row = 45; Price = "0.593";      Volume = "  +0.13%  " },
    @{ Row = 46; Price = $null;        Volume = "  -0.31%  " },
    @{ Row = 47; Price = "0.0532";     Volume = "  +0.96%  " },
    @{ Row = 48; Price = "18.99";      Volume = "  +2.21%  " },
    @{ Row = 49; Price = "123.91";     Volume = "  +11.06%  " },
    @{ Row = 50; Price = "0.0231";     Volume = "  +1.95%  " },
    @{ Row = 51; Price = "18.48";      Volume = "  +4.07%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Cells.Item($u.Row, 4)   # column D = Price
        $text = $u.Price
        $looksNumeric = $text -match '^[0-9]+(\.[0-9]+)?$'
        if ($looksNumeric) {
            # Force text storage (leading apostrophe / quote-prefix) so the
            # value isn't silently re-typed as a Number by Excel's
            # auto-detect, keeping parity with the rest of the text column.
            $priceCell.Value = "'" + $text
        } else {
            $priceCell.Value = $text
        }
    }

    $volumeCell = $ws.Cells.Item($u.Row, 5)      # column E = Volume(1h)
    $volumeCell.Value = $u.Volume
}
